$d = $word.ActiveDocument

# 1. Insert a new run containing "hell" immediately before the existing
#    "kdjfkfjk" run at the very start of the Title paragraph.
$firstPara = $d.Paragraphs(1)
$startRng = $firstPara.Range
$startRng.SetRange(0, 0)
$startRng.InsertBefore("hell")

# 2. Append a brand-new, plain (unstyled) paragraph "My name is yasir" at
#    the end of the document body, after the existing paragraph and
#    before the sectPr. We build it from raw OOXML via InsertXML so it
#    does not inherit the Title style / superscript run formatting that
#    sits at the end of the last paragraph.
$docEnd = $d.Content.End
$insertionPoint = $d.Range($docEnd, $docEnd)
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>My name is yasir</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($newParaXml)
